# Ejecución y ajuste flujo de aplicación de pagos - 10 Pagadurias
# Updates the single data row in the (only) worksheet to reflect the
# new "pagaduria" values used for this test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 = FechaRegistro : 10/12/2021 -> 22/12/2021
$ws.Range("J2").Value = "22/12/2021"

# A2 = IdPagaduria : 103 -> 349
$ws.Range("A2").Value = "349"

# C2 = NombrePagaduria : "CONSORCIO DE PENSIONES DEL HUILA" ->
#      "ALCALDÍA MUNICIPAL DE BARRANCABERMEJA NÓMINA TRABAJADORES OFICIALES"
$ws.Range("C2").Value = '"ALCALDÍA MUNICIPAL DE BARRANCABERMEJA NÓMINA TRABAJADORES OFICIALES"'

# Move the active selection to C7, matching the author's saved cursor position.
$ws.Range("C7").Select()
